# feat: add 2022-Q1 data
#
# Inserts a new "2022-Q1" worksheet (between "2021-Q3" and "总计") with per-
# fund holding detail, and records the quarter's summary row at the top of
# the "总计" sheet.

$wb = $excel.ActiveWorkbook

$wsQ3    = $wb.Worksheets.Item(1)
$wsTotal = $wb.Worksheets.Item(2)

# --- 1. Add the new "2022-Q1" sheet right after "2021-Q3" -------------------
# NOTE: once a sheet is inserted, any previously-held worksheet variable that
# resolves by position can end up pointing at the newly inserted sheet, so we
# re-fetch "总计" by name after this structural change instead of reusing
# $wsTotal.
$wsQ1 = $wb.Worksheets.Add($wsTotal)
$wsQ1.Name = "2022-Q1"

# Header row (B1:H1).
$wsQ1.Range("B1").Value = "基金代码"
$wsQ1.Range("C1").Value = "基金名称"
$wsQ1.Range("D1").Value = "基金规模"
$wsQ1.Range("E1").Value = "股票总仓位"
$wsQ1.Range("F1").Value = "仓位占比"
$wsQ1.Range("G1").Value = "持有市值(亿元)"
$wsQ1.Range("H1").Value = "仓位排名"

# Row 2 - 招商量化精选股票A
$wsQ1.Range("A2").Value = 0
$wsQ1.Range("B2").Value = "'001917"
$wsQ1.Range("B2").Style = "Normal"
$wsQ1.Range("C2").Value = "招商量化精选股票A"
$wsQ1.Range("D2").Value = "'2.33"
$wsQ1.Range("D2").Style = "Normal"
$wsQ1.Range("E2").Value = "'94.20"
$wsQ1.Range("E2").Style = "Normal"
$wsQ1.Range("F2").Value = "'1.23"
$wsQ1.Range("F2").Style = "Normal"
$wsQ1.Range("G2").Value = "'0.0287"
$wsQ1.Range("G2").Style = "Normal"
$wsQ1.Range("H2").Value = 7

# Row 3 - 招商量化精选股票C
$wsQ1.Range("A3").Value = 1
$wsQ1.Range("B3").Value = "'007950"
$wsQ1.Range("B3").Style = "Normal"
$wsQ1.Range("C3").Value = "招商量化精选股票C"
$wsQ1.Range("D3").Value = "'0.56"
$wsQ1.Range("D3").Style = "Normal"
$wsQ1.Range("E3").Value = "'94.20"
$wsQ1.Range("E3").Style = "Normal"
$wsQ1.Range("F3").Value = "'1.23"
$wsQ1.Range("F3").Style = "Normal"
$wsQ1.Range("G3").Value = "'0.0069"
$wsQ1.Range("G3").Style = "Normal"
$wsQ1.Range("H3").Value = 7

# Copy the "总计" header formatting (bold + border, style used on B1/A2 of
# that sheet) onto the header row and the index column of the new sheet.
$wsTotalFmt = $wb.Worksheets.Item("总计")
$wsTotalFmt.Range("B1").Copy()
$wsQ1.Range("B1:H1").PasteSpecial(-4122)
$wsTotalFmt.Range("A2").Copy()
$wsQ1.Range("A2:A3").PasteSpecial(-4122)

# --- 2. Insert the 2022-Q1 summary row at the top of "总计" -----------------
$wsTotal2 = $wb.Worksheets.Item("总计")
$wsTotal2.Rows.Item(2).Insert()

# Inserting a row re-uses nearby formatting for the blank cells it creates
# (e.g. it picks up the header row's bold/centered look for B2:D2) - reset
# that back to the default "Normal" style before writing the real values, to
# match the unstyled data cells used elsewhere on this sheet.
$wsTotal2.Range("B2:D2").Style = "Normal"

$wsTotal2.Range("B2").Value = "2022-Q1"
$wsTotal2.Range("C2").Value = 2
$wsTotal2.Range("D2").Value = 0.04

# A2/A3 are the running row index (0, 1, 2, ...) - renumber them and restore
# the bordered index-column style (copied from the row that now holds the
# old 2021-Q3 data, i.e. row 3) onto the newly inserted row 2.
$wsTotal2.Range("A3").Copy()
$wsTotal2.Range("A2").PasteSpecial(-4122)
$wsTotal2.Range("A2").Value = 0
$wsTotal2.Range("A3").Value = 1
